$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-21 Tuesday", 2)

$d.Content.Find.Execute("715÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "452÷8=", 2)
$d.Content.Find.Execute("645÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "105÷6=", 2)
$d.Content.Find.Execute("638÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "148÷8=", 2)
$d.Content.Find.Execute("203÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "885÷8=", 2)
$d.Content.Find.Execute("737÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "643÷2=", 2)

$d.Content.Find.Execute("419÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "112÷7=", 2)
$d.Content.Find.Execute("866÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "584÷6=", 2)
$d.Content.Find.Execute("384÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "309÷6=", 2)
$d.Content.Find.Execute("573÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "262÷6=", 2)
$d.Content.Find.Execute("653÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "583÷5=", 2)

$d.Content.Find.Execute("511÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "857÷6=", 2)
$d.Content.Find.Execute("635÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "138÷4=", 2)
$d.Content.Find.Execute("964÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "138÷3=", 2)
$d.Content.Find.Execute("355÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "374÷4=", 2)
$d.Content.Find.Execute("871÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "215÷5=", 2)

$d.Content.Find.Execute("651÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "321÷4=", 2)
$d.Content.Find.Execute("412÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "136÷3=", 2)
$d.Content.Find.Execute("664÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "157÷3=", 2)
$d.Content.Find.Execute("293÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "738÷7=", 2)
$d.Content.Find.Execute("670÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "387÷6=", 2)

$d.Content.Find.Execute("424÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "361÷7=", 2)
$d.Content.Find.Execute("323÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "264÷5=", 2)
$d.Content.Find.Execute("704÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "422÷4=", 2)
$d.Content.Find.Execute("608÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "855÷6=", 2)
$d.Content.Find.Execute("560÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "311÷3=", 2)
